$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.568.94"
$ws.Range("E2").Value = "  +4.59%  "
$ws.Range("D3").Value = "2.349.90"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "546.89"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").Value = "132.50"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("D9").Value = "2.346.86"
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "0.334"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("E14").Value = "  +2.50%  "
$ws.Range("D15").Value = "2.766.84"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").Value = "60.517.46"
$ws.Range("E16").Value = "  +4.62%  "
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("D18").Value = "2.343.86"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").Value = "10.69"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").Value = "6.89"
$ws.Range("E21").Value = "  +8.54%  "
$ws.Range("D22").Value = "314.67"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "63.33"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("E25").Value = "  +3.23%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "7.94"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").Value = "1.35"
$ws.Range("E28").Value = "  +5.29%  "
$ws.Range("E29").Value = "  +2.85%  "
$ws.Range("D30").Value = "171.65"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("E31").Value = "  +10.12%  "
$ws.Range("E32").Value = "  +2.12%  "
$ws.Range("D33").Value = "5.91"
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("D34").Value = "1.41"
$ws.Range("E34").Value = "  +14.97%  "
$ws.Range("D35").Value = "0.381"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").Value = "18.06"
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +7.19%  "
$ws.Range("D40").Value = "314.41"
$ws.Range("E40").Value = "  +9.60%  "
$ws.Range("D41").Value = "38.19"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  +3.48%  "
$ws.Range("D43").Value = "142.50"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").Value = "3.47"
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("D45").Value = "0.0956"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").Value = "19.29"
$ws.Range("E46").Value = "  +6.85%  "
$ws.Range("D47").Value = "0.0498"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").Value = "0.561"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").Value = "0.0₆0208"
$ws.Range("E51").Value = "  +4.01%  "
